# Apply the weekly-programs update: new date range (05/03/2024-30/04/2024),
# refreshed Bible-reading / song / talk assignments, and two additional
# weeks appended (rows 16-18), shifting the closing rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 holds dd/mm/yyyy text dates; force Text format first so Excel
# does not reinterpret ambiguous values (e.g. 05/03/2024) as mm/dd dates.
$ws.Range("A1:I1").NumberFormat = "@"

# Row 1
$ws.Range("A1").Value = "05/03/2024"
$ws.Range("B1").Value = "12/03/2024"
$ws.Range("C1").Value = "19/03/2024"
$ws.Range("D1").Value = "26/03/2024"
$ws.Range("E1").Value = "02/04/2024"
$ws.Range("F1").Value = "09/04/2024"
$ws.Range("G1").Value = "16/04/2024"
$ws.Range("H1").Value = "23/04/2024"
$ws.Range("I1").Value = "30/04/2024"
# Row 2
$ws.Range("A2").Value = "SALMOS 16,17"
$ws.Range("B2").Value = "SALMO 18"
$ws.Range("C2").Value = "SALMOS 19-21"
$ws.Range("D2").Value = "SALMO 22"
$ws.Range("E2").Value = "SALMOS 23-25"
$ws.Range("F2").Value = "SALMOS 26-28"
$ws.Range("G2").Value = "SALMOS 29-31"
$ws.Range("H2").Value = "SALMOS 32,33"
$ws.Range("I2").Value = "SALMOS 34,35"
# Row 3
$ws.Range("A3").Value = "Canción 111y oración | Palabras de introducción(1 min.)"
$ws.Range("B3").Value = "Canción 148y oración | Palabras de introducción(1 min.)"
$ws.Range("C3").Value = "Canción 6y oración | Palabras de introducción(1 min.)"
$ws.Range("D3").Value = "Canción 19y oración | Palabras de introducción(1 min.)"
$ws.Range("E3").Value = "Canción 4y oración | Palabras de introducción(1 min.)"
$ws.Range("F3").Value = "Canción 34y oración | Palabras de introducción(1 min.)"
$ws.Range("G3").Value = "Canción 108y oración | Palabras de introducción(1 min.)"
$ws.Range("H3").Value = "Canción 103y oración | Palabras de introducción(1 min.)"
$ws.Range("I3").Value = "Canción 10y oración | Palabras de introducción(1 min.)"
# Row 4
$ws.Range("A4").Value = "TESOROS DE LA BIBLIA"
$ws.Range("B4").Value = "TESOROS DE LA BIBLIA"
$ws.Range("C4").Value = "TESOROS DE LA BIBLIA"
$ws.Range("D4").Value = "TESOROS DE LA BIBLIA"
$ws.Range("E4").Value = "TESOROS DE LA BIBLIA"
$ws.Range("F4").Value = "TESOROS DE LA BIBLIA"
$ws.Range("G4").Value = "TESOROS DE LA BIBLIA"
$ws.Range("H4").Value = "TESOROS DE LA BIBLIA"
$ws.Range("I4").Value = "TESOROS DE LA BIBLIA"
# Row 5
$ws.Range("A5").Value = "1. “Jehová, mi fuente de todo lo bueno”"
$ws.Range("B5").Value = "1. “Jehová es [...] mi libertador”"
$ws.Range("C5").Value = "1. “Los cielos declaran la gloria de Dios”"
$ws.Range("D5").Value = "1. Se profetizaron detalles de la muerte de Jesús"
$ws.Range("E5").Value = "1. “Jehová es mi Pastor”"
$ws.Range("F5").Value = "1. Qué ayudaba a David a vivir con integridad"
$ws.Range("G5").Value = "1. La disciplina es una muestra del amor de Dios"
$ws.Range("H5").Value = "1. ¿Por qué deben confesarse los pecados graves?"
$ws.Range("I5").Value = "1. “Alabaré a Jehová en todo momento”"
# Row 6
$ws.Range("A6").Value = "2. Busquemos perlas escondidas"
$ws.Range("B6").Value = "2. Busquemos perlas escondidas"
$ws.Range("C6").Value = "2. Busquemos perlas escondidas"
$ws.Range("D6").Value = "2. Busquemos perlas escondidas"
$ws.Range("E6").Value = "2. Busquemos perlas escondidas"
$ws.Range("F6").Value = "2. Busquemos perlas escondidas"
$ws.Range("G6").Value = "2. Busquemos perlas escondidas"
$ws.Range("H6").Value = "2. Busquemos perlas escondidas"
$ws.Range("I6").Value = "2. Busquemos perlas escondidas"
# Row 7
$ws.Range("A7").Value = "3. Lectura de la Biblia"
$ws.Range("B7").Value = "3. Lectura de la Biblia"
$ws.Range("C7").Value = "3. Lectura de la Biblia"
$ws.Range("D7").Value = "3. Lectura de la Biblia"
$ws.Range("E7").Value = "3. Lectura de la Biblia"
$ws.Range("F7").Value = "3. Lectura de la Biblia"
$ws.Range("G7").Value = "3. Lectura de la Biblia"
$ws.Range("H7").Value = "3. Lectura de la Biblia"
$ws.Range("I7").Value = "3. Lectura de la Biblia"
# Row 8
$ws.Range("A8").Value = "SEAMOS MEJORES MAESTROS"
$ws.Range("B8").Value = "SEAMOS MEJORES MAESTROS"
$ws.Range("C8").Value = "SEAMOS MEJORES MAESTROS"
$ws.Range("D8").Value = "SEAMOS MEJORES MAESTROS"
$ws.Range("E8").Value = "SEAMOS MEJORES MAESTROS"
$ws.Range("F8").Value = "SEAMOS MEJORES MAESTROS"
$ws.Range("G8").Value = "SEAMOS MEJORES MAESTROS"
$ws.Range("H8").Value = "SEAMOS MEJORES MAESTROS"
$ws.Range("I8").Value = "SEAMOS MEJORES MAESTROS"
# Row 9
$ws.Range("A9").Value = "4. Empiece conversaciones"
$ws.Range("B9").Value = "4. Bondad: Lo que hizo Jesús"
$ws.Range("C9").Value = "4. Empiece conversaciones"
$ws.Range("D9").Value = "4. Empiece conversaciones"
$ws.Range("E9").Value = "4. Empiece conversaciones"
$ws.Range("F9").Value = "4. Empiece conversaciones"
$ws.Range("G9").Value = "4. Empiece conversaciones"
$ws.Range("H9").Value = "4. Humildad: Lo que hizo Pablo"
$ws.Range("I9").Value = "4. Empiece conversaciones"
# Row 10
$ws.Range("A10").Value = "5. Empiece conversaciones"
$ws.Range("B10").Value = "5. Bondad: Imite a Jesús"
$ws.Range("C10").Value = "5. Empiece conversaciones"
$ws.Range("D10").Value = "5. Haga revisitas"
$ws.Range("E10").Value = "5. Haga revisitas"
$ws.Range("F10").Value = "5. Haga revisitas"
$ws.Range("G10").Value = "5. Empiece conversaciones"
$ws.Range("H10").Value = "5. Humildad: Imite a Pablo"
$ws.Range("I10").Value = "5. Haga revisitas"
# Row 11
$ws.Range("A11").Value = "6. Empiece conversaciones"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "6. Explique sus creencias"
$ws.Range("D11").Value = "6. Discurso"
$ws.Range("E11").Value = "6. Haga discípulos"
$ws.Range("F11").Value = "6. Discurso"
$ws.Range("G11").Value = "6. Haga revisitas"
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = "6. Explique sus creencias"
# Row 12
$ws.Range("A12").Value = "7. Haga discípulos"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = "7. Haga discípulos"
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""
# Row 13
$ws.Range("A13").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("B13").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("C13").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("D13").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("E13").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("F13").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("G13").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("H13").Value = "NUESTRA VIDA CRISTIANA"
$ws.Range("I13").Value = "NUESTRA VIDA CRISTIANA"
# Row 14
$ws.Range("A14").Value = "Canción 20"
$ws.Range("B14").Value = "Canción 60"
$ws.Range("C14").Value = "Canción 141"
$ws.Range("D14").Value = "Canción 95"
$ws.Range("E14").Value = "Canción 54"
$ws.Range("F14").Value = "Canción 128"
$ws.Range("G14").Value = "Canción 45"
$ws.Range("H14").Value = "Canción 74"
$ws.Range("I14").Value = "Canción 59"
# Row 15
$ws.Range("A15").Value = "8. ¡Preparémonos para la Conmemoración!"
$ws.Range("B15").Value = "6. Necesidades de la congregación"
$ws.Range("C15").Value = "7.Observar la creación fortalece la fe"
$ws.Range("D15").Value = "7. Necesidades de la congregación"
$ws.Range("E15").Value = "7. Rechazamos la voz de los extraños"
$ws.Range("F15").Value = "7. Adolescentes que son moralmente íntegros"
$ws.Range("G15").Value = "8.Por qué tenemos fe en... el amor de Dios"
$ws.Range("H15").Value = "6. Necesidades de la congregación"
$ws.Range("I15").Value = "7. Tres formas de alabar a Jehová en nuestras reuniones"
# Row 16
$ws.Range("A16").Value = "9. Estudio bíblico de la congregación"
$ws.Range("B16").Value = "7.Logros de la organizaciónpara el mes de marzo"
$ws.Range("C16").Value = "8. Estudio bíblico de la congregación"
$ws.Range("D16").Value = "8. Estudio bíblico de la congregación"
$ws.Range("E16").Value = "8. Estudio bíblico de la congregación"
$ws.Range("F16").Value = "8. Estudio bíblico de la congregación"
$ws.Range("G16").Value = "9.2024|Informe sobre la actividad del Departamento Local de Diseño y Construcción"
$ws.Range("H16").Value = "7. Estudio bíblico de la congregación"
$ws.Range("I16").Value = "8. Estudio bíblico de la congregación"
# Row 17
$ws.Range("A17").Value = ""
$ws.Range("B17").Value = "8. Estudio bíblico de la congregación"
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("G17").Value = "10. Estudio bíblico de la congregación"
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""
# Row 18
$ws.Range("A18").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""

Write-Host "Applied weekly_programs update"
